# Remove form_id from remaining forms
#
# The "settings" sheet has a "form_id" column (column B) that is no longer
# needed; this script removes it entirely. Removing the column shifts the
# "version", "style" and "namespaces" columns one place to the left (C->B,
# D->C, E->D), along with their cell-comment annotations, which we migrate
# by hand since cell comments stay anchored to their original address when
# a column is deleted.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

# 1. Shift the header-row comments left by one, matching the columns that
#    are about to shift left when form_id (column B) is deleted.
$versionComment    = $ws.Range("C1").Comment.Text()
$pagesComment      = $ws.Range("D1").Comment.Text()
$namespacesComment = $ws.Range("E1").Comment.Text()

[void]$ws.Range("B1").Comment.Text($versionComment)
[void]$ws.Range("C1").Comment.Text($pagesComment)
[void]$ws.Range("D1").Comment.Text($namespacesComment)

# The old E1 comment (namespaces) has now been copied onto D1, so the
# original is redundant and can be dropped.
$ws.Range("E1").Comment.Delete()

# 2. Delete the form_id column itself (column B). This shifts the
#    version/style/namespaces cell values and the formula cell left by one
#    column, same as the comment text above.
$ws.Columns.Item(2).Delete()
